# Auto-generated Excel COM-interop script to apply the scraped schedule update
# (Linea 141 horarios refresh: new scrape timestamp 12:44:21, 14 new rows on LP1912,
# 1 new row on LP1912-215, and reordered/updated data rows reflecting the new scrape).

$wb = $excel.ActiveWorkbook

# ---- LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = 'Última actualización: 12:44:21'
$ws1.Range("A3").Value = 'Total filas: 260'
$ws1.Range("C52").Value = '16_SANTA ANA'
$ws1.Range("C53").Value = '11_ETCHEVERRY'
$ws1.Range("C118").Value = '16_SANTA ANA'
$ws1.Range("C119").Value = '23_HERNANDEZ'
$ws1.Range("C199").Value = '14_ABASTO'
$ws1.Range("C201").Value = '16_P MOR-SANTA ANA'
$ws1.Range("A208").Value = '10:28:12'
$ws1.Range("C208").Value = '17_ROMERO'
$ws1.Range("D208").Value = 106
$ws1.Range("A209").Value = '10:57:58'
$ws1.Range("C209").Value = '10_OLMOS'
$ws1.Range("D209").Value = 77
$ws1.Range("C214").Value = '215A_EL PATO'
$ws1.Range("A215").Value = '10:28:12'
$ws1.Range("C215").Value = '26_HERNANDEZ'
$ws1.Range("D215").Value = 113
$ws1.Range("A216").Value = '11:51:05'
$ws1.Range("C216").Value = '14_ABASTO'
$ws1.Range("D216").Value = 30
$ws1.Range("A217").Value = '12:16:51'
$ws1.Range("C217").Value = '16_SANTA ANA'
$ws1.Range("D217").Value = 5
$ws1.Range("A230").Value = '12:44:21'
$ws1.Range("B230").Value = '12:44'
$ws1.Range("C230").Value = '16_SANTA ANA'
$ws1.Range("D230").Value = 0
$ws1.Range("A231").Value = '12:44:21'
$ws1.Range("B231").Value = '12:44'
$ws1.Range("C231").Value = '10_OLMOS'
$ws1.Range("D231").Value = 0
$ws1.Range("A232").Value = '12:44:21'
$ws1.Range("B232").Value = '12:45'
$ws1.Range("C232").Value = '11_ETCHEVERRY'
$ws1.Range("D232").Value = 1
$ws1.Range("A233").Value = '12:44:21'
$ws1.Range("B233").Value = '12:47'
$ws1.Range("C233").Value = '16_SANTA ANA'
$ws1.Range("D233").Value = 3
$ws1.Range("B234").Value = '12:48'
$ws1.Range("C234").Value = '11_ETCHEVERRY'
$ws1.Range("D234").Value = 85
$ws1.Range("B235").Value = '12:49'
$ws1.Range("C235").Value = '11_ETCHEVERRY'
$ws1.Range("D235").Value = 58
$ws1.Range("A236").Value = '11:23:54'
$ws1.Range("B236").Value = '12:54'
$ws1.Range("C236").Value = '17_ROMERO'
$ws1.Range("D236").Value = 91
$ws1.Range("A237").Value = '12:44:21'
$ws1.Range("B237").Value = '13:02'
$ws1.Range("C237").Value = '14_ABASTO'
$ws1.Range("D237").Value = 18
$ws1.Range("B238").Value = '13:02'
$ws1.Range("C238").Value = '15_ABASTO'
$ws1.Range("D238").Value = 71
$ws1.Range("B239").Value = '13:06'
$ws1.Range("C239").Value = '16_P MOR-SANTA ANA'
$ws1.Range("D239").Value = 103
$ws1.Range("A240").Value = '11:51:05'
$ws1.Range("B240").Value = '13:07'
$ws1.Range("C240").Value = '16_P MOR-SANTA ANA'
$ws1.Range("D240").Value = 76
$ws1.Range("A241").Value = '12:16:51'
$ws1.Range("B241").Value = '13:08'
$ws1.Range("D241").Value = 52
$ws1.Range("A242").Value = '11:23:54'
$ws1.Range("B242").Value = '13:13'
$ws1.Range("C242").Value = '215D_EL PATO'
$ws1.Range("D242").Value = 110
$ws1.Range("A243").Value = '12:44:21'
$ws1.Range("B243").Value = '13:14'
$ws1.Range("C243").Value = '11_ETCHEVERRY'
$ws1.Range("D243").Value = 30
$ws1.Range("B244").Value = '13:14'
$ws1.Range("C244").Value = '215D_EL PATO'
$ws1.Range("D244").Value = 83
$ws1.Range("A245").Value = '11:23:54'
$ws1.Range("B245").Value = '13:19'
$ws1.Range("D245").Value = 116
$ws1.Range("A246").Value = '11:51:05'
$ws1.Range("B246").Value = '13:20'
$ws1.Range("C246").Value = '10_OLMOS'
$ws1.Range("D246").Value = 89
$ws1.Range("A247").Value = '11:23:54'
$ws1.Range("B247").Value = '13:20'
$ws1.Range("C247").Value = '26_HERNANDEZ'
$ws1.Range("D247").Value = 117
$ws1.Range("A248").Value = '12:44:21'
$ws1.Range("B248").Value = '13:21'
$ws1.Range("C248").Value = '10_OLMOS'
$ws1.Range("D248").Value = 37
$ws1.Range("A249").Value = '11:51:05'
$ws1.Range("B249").Value = '13:21'
$ws1.Range("C249").Value = '26_HERNANDEZ'
$ws1.Range("D249").Value = 90
$ws1.Range("B250").Value = '13:26'
$ws1.Range("C250").Value = '14_ABASTO'
$ws1.Range("D250").Value = 70
$ws1.Range("A251").Value = '11:51:05'
$ws1.Range("B251").Value = '13:27'
$ws1.Range("C251").Value = '14_ABASTO'
$ws1.Range("D251").Value = 96
$ws1.Range("A252").Value = '12:16:51'
$ws1.Range("B252").Value = '13:32'
$ws1.Range("C252").Value = '10_OLMOS'
$ws1.Range("D252").Value = 76
$ws1.Range("E252").Value = 'LP1912'
$ws1.Range("A253").Value = '12:16:51'
$ws1.Range("B253").Value = '13:34'
$ws1.Range("C253").Value = '23_HERNANDEZ'
$ws1.Range("D253").Value = 78
$ws1.Range("E253").Value = 'LP1912'
$ws1.Range("A254").Value = '12:44:21'
$ws1.Range("B254").Value = '13:35'
$ws1.Range("C254").Value = '23_HERNANDEZ'
$ws1.Range("D254").Value = 51
$ws1.Range("E254").Value = 'LP1912'
$ws1.Range("A255").Value = '11:51:05'
$ws1.Range("B255").Value = '13:36'
$ws1.Range("C255").Value = '15_ABASTO'
$ws1.Range("D255").Value = 105
$ws1.Range("E255").Value = 'LP1912'
$ws1.Range("A256").Value = '11:51:05'
$ws1.Range("B256").Value = '13:46'
$ws1.Range("C256").Value = '17_ROMERO'
$ws1.Range("D256").Value = 115
$ws1.Range("E256").Value = 'LP1912'
$ws1.Range("A257").Value = '12:16:51'
$ws1.Range("B257").Value = '13:50'
$ws1.Range("C257").Value = '215A_EL PATO'
$ws1.Range("D257").Value = 94
$ws1.Range("E257").Value = 'LP1912'
$ws1.Range("A258").Value = '12:16:51'
$ws1.Range("B258").Value = '13:55'
$ws1.Range("C258").Value = '225_GOMEZ'
$ws1.Range("D258").Value = 99
$ws1.Range("E258").Value = 'LP1912'
$ws1.Range("A259").Value = '12:44:21'
$ws1.Range("B259").Value = '13:56'
$ws1.Range("C259").Value = '225_GOMEZ'
$ws1.Range("D259").Value = 72
$ws1.Range("E259").Value = 'LP1912'
$ws1.Range("A260").Value = '12:16:51'
$ws1.Range("B260").Value = '14:04'
$ws1.Range("C260").Value = '17_ROMERO'
$ws1.Range("D260").Value = 108
$ws1.Range("E260").Value = 'LP1912'
$ws1.Range("A261").Value = '12:44:21'
$ws1.Range("B261").Value = '14:05'
$ws1.Range("C261").Value = '23_HERNANDEZ'
$ws1.Range("D261").Value = 81
$ws1.Range("E261").Value = 'LP1912'
$ws1.Range("A262").Value = '12:44:21'
$ws1.Range("B262").Value = '14:13'
$ws1.Range("C262").Value = '16_P MOR-167 Y 521'
$ws1.Range("D262").Value = 89
$ws1.Range("E262").Value = 'LP1912'
$ws1.Range("A263").Value = '12:44:21'
$ws1.Range("B263").Value = '14:17'
$ws1.Range("C263").Value = '27_EL RETIRO'
$ws1.Range("D263").Value = 93
$ws1.Range("E263").Value = 'LP1912'
$ws1.Range("A264").Value = '12:44:21'
$ws1.Range("B264").Value = '14:20'
$ws1.Range("C264").Value = '215C_EL PATO'
$ws1.Range("D264").Value = 96
$ws1.Range("E264").Value = 'LP1912'
$ws1.Range("A265").Value = '12:44:21'
$ws1.Range("B265").Value = '14:21'
$ws1.Range("C265").Value = '26_HERNANDEZ'
$ws1.Range("D265").Value = 97
$ws1.Range("E265").Value = 'LP1912'

# ---- LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = 'Última actualización: 12:44:21'
$ws2.Range("A3").Value = 'Total filas: 29'
$ws2.Range("A34").Value = '12:44:21'
$ws2.Range("B34").Value = '14:20'
$ws2.Range("C34").Value = '215C_EL PATO'
$ws2.Range("D34").Value = 96
$ws2.Range("E34").Value = 'LP1912'

# ---- 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = 'Última actualización: 12:44:21'

